$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the input parameter values (D2:G3) per the new "NN Size" results
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 2

$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

# Update the selected cell shown in the sheet view
$ws.Range("B6").Select()
